$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "parameters" (sheet2.xml)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("parameters")

# Insert two new rows right before the existing "Maintance & repair cost" row
# (row 25), pushing it and everything below it down by two rows. Excel will
# auto-adjust the =123500/E27 formula in E2 to follow the "Cost to price
# markup factor" row as it shifts from row 27 to row 29.
$ws.Rows("25:26").Insert()

# New row 25: Dwell time basis
$ws.Range("A25").Value = "Class 8 Diesel Tractor"
$ws.Range("B25").Value = "Reference"
$ws.Range("C25").Value = "Dwell time basis"
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = "n/a"
$ws.Range("G25").Value = "Selection of fuel/energy storage for dwell time calculation (0 = fuel storage, 1 = battery)"

# New row 26: Dwell time boolean
$ws.Range("A26").Value = "Class 8 Diesel Tractor"
$ws.Range("B26").Value = "Reference"
$ws.Range("C26").Value = "Dwell time boolean"
$ws.Range("D26").Value = 24
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = "n/a"
$ws.Range("G26").Value = "Boolean to turn on or off dwell time costs (1 = on, 0 = off)"

# Append two brand-new rows (32 and 33) after the existing data (which now
# runs through row 31 post-insert).
$ws.Range("A32").Value = "Class 8 Diesel Tractor"
$ws.Range("B32").Value = "Reference"
$ws.Range("C32").Value = "Carbon cost boolean"
$ws.Range("D32").Value = 30
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = "n/a"
$ws.Range("G32").Value = "Boolean to turn on or off carbon costs (1 = on, 0 = off)"

$ws.Range("A33").Value = "Class 8 Diesel Tractor"
$ws.Range("B33").Value = "Reference"
$ws.Range("C33").Value = "Input fuel efficiency"
$ws.Range("D33").Value = 31
$ws.Range("E33").Formula = "=1/(7*1.136)"
$ws.Range("F33").Value = "gge/mile"
$ws.Range("G33").Value = "Inverse of diesel fuel economy of 7 mile/dge (duplicates design variable)"

$ws.Range("C34").Select()

# ---------------------------------------------------------------------------
# Sheet "designs" (sheet1.xml)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("designs")
$ws1.Columns("C").ColumnWidth = 15
$ws1.Range("G14").Select()

# ---------------------------------------------------------------------------
# Sheet "results" (sheet3.xml) - swap two pairs of metric rows
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("results")

$c4 = $ws3.Range("C4").Value2
$d4 = $ws3.Range("D4").Value2
$c5 = $ws3.Range("C5").Value2
$d5 = $ws3.Range("D5").Value2
$ws3.Range("C4").Value = $c5
$ws3.Range("D4").Value = $d5
$ws3.Range("C5").Value = $c4
$ws3.Range("D5").Value = $d4

$c6 = $ws3.Range("C6").Value2
$d6 = $ws3.Range("D6").Value2
$c7 = $ws3.Range("C7").Value2
$d7 = $ws3.Range("D7").Value2
$ws3.Range("C6").Value = $c7
$ws3.Range("D6").Value = $d7
$ws3.Range("C7").Value = $c6
$ws3.Range("D7").Value = $d6

$ws3.Range("C3:C7").Select()

# ---------------------------------------------------------------------------
# Sheet "indices" (sheet4.xml) - swap two pairs of Metric labels (index
# numbers in column D stay put; only the label in column C moves)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("indices")

$c12 = $ws4.Range("C12").Value2
$c13 = $ws4.Range("C13").Value2
$ws4.Range("C12").Value = $c13
$ws4.Range("C13").Value = $c12

$c14 = $ws4.Range("C14").Value2
$c15 = $ws4.Range("C15").Value2
$ws4.Range("C14").Value = $c15
$ws4.Range("C15").Value = $c14

$ws4.Range("C11").Select()
